$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Date" property value (row 8, column B) - refreshed generation timestamp
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"

# "Jurisdiction" property value (row 11, column B) - was blank, now set to FRANCE
$ws.Range("B11").Value = "FRANCE"
